# Weekly fruit/vegetable price data refresh.
# The underlying daily records (columns D, J, K, L, M, O, P) were reshuffled
# across rows 2-16 (rows 8 and 11 stay as-is). Columns A,B,C,E,F,G,H,I,N,Q,R
# are identical for every data row, so only the changed columns are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target row number -> @(D, J, K, L, M, O, P)
$rows = @{
    2  = @(44650, 130, 3000, 3500, 3308, "Región Metropolitana", 551)
    3  = @(44685, 150, 3000, 3500, 3267, "Región Metropolitana", 544)
    4  = @(44672, 140, 3000, 3500, 3286, "Región Metropolitana", 548)
    5  = @(44658, 180, 2500, 3000, 2778, "Región Metropolitana", 463)
    6  = @(44876,  80, 6500, 7000, 6812, "Región Metropolitana", 1135)
    7  = @(44671, 150, 3500, 4000, 3733, "Región Metropolitana", 622)
    9  = @(44659,  90, 2500, 3000, 2722, "Región Metropolitana", 454)
    10 = @(44643,  90, 2800, 3000, 2911, "Región Metropolitana", 485)
    12 = @(44631, 110, 3000, 3500, 3273, "Provincia de Chacabuco", 546)
    13 = @(44644, 140, 2500, 3000, 2786, "Provincia de Chacabuco", 464)
    14 = @(44630,  90, 2500, 3000, 2722, "Región Metropolitana", 454)
    15 = @(44637, 170, 2800, 3000, 2906, "Región Metropolitana", 484)
    16 = @(44957,  70, 1500, 2000, 1857, "Región Metropolitana", 310)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $vals[5]   # O - Origen
    $ws.Cells.Item($r, 16).Value = $vals[6]   # P - Precio $/Kg
}
